# Apply the "Changed type kind values" edit to the "meta types" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta types")

# Rows 22-37: bump the Value column (B) by 10.
for ($r = 22; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value2 + 10
}

# Rows 38-62: bump the Value column (B) by 30.
for ($r = 38; $r -le 62; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value2 + 30
}

# Rows 36 and 37 (tkMetaType / tkMetaRepo) no longer list the C++ type
# (metapp::MetaType / metapp::MetaRepo) - remove those cells entirely.
$ws.Range("C36").ClearContents()
$ws.Range("C37").ClearContents()

# Update the saved view state: scrolled down so row 29 is at the top,
# with C37 selected (was B18 at the top with B18 selected).
$ws.Range("C37").Select()
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 29
$activeWindow.ScrollColumn = 1
